# Deploy_etb_idartes_18_noviembre_2025 - "Entrega 18 noviembre 2025"
#
# This script reproduces, via Word COM-interop calls, the edits described
# by the target diff:
#   1) Splits the "Levantar servicios:" bullet into three paragraphs,
#      inserting a new bold "Permisos de carpeta uploads:" bullet and a
#      "* Ejecutar: chmod -R 775 v1/uploads/" sub-paragraph (carrying the
#      "_GoBack" bookmark) right before it.
#   2) Adds a <w:lastRenderedPageBreak/> before "PROJECT_VERSION=v1".
#   3) Adds a <w:lastRenderedPageBreak/> before
#      "AWS_ACCESS_KEY_ID=AKIAXYKJTSSH3VRHA7HC".
#   4) Removes the old "_GoBack" bookmark that used to sit after
#      "... i dentro de la carpeta v1" (it moved to the new location
#      created in step 1, since Word keeps a single "_GoBack" bookmark
#      tracking the most recent edit point).

$d = $word.ActiveDocument

function Get-ParagraphByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    throw "Paragraph containing '$needle' not found"
}

# ---------------------------------------------------------------------
# 1) "Levantar servicios:" -> split into three paragraphs
# ---------------------------------------------------------------------
$levantar = Get-ParagraphByText("Levantar servicios")
$levantarXml = (
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Permisos de carpeta </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>uploads</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:br/></w:r>' +
    '<w:r><w:t>Se debe salvar para futuros despliegues con el objeto de no perder datos y garantizar permisos totales sobre la carpeta v1/</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>uploads</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">/ que sería la </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>storage</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> del proyecto</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="720"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>*</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Ejecutar</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>chmod</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> -R 77</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>5 v1/uploads/</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Levantar servicios: </w:t></w:r>' +
    '<w:r><w:t>Ejecu</w:t></w:r>' +
    '<w:r><w:t>tar según configuración en servidor o por ambiente bajo los comandos especificados.</w:t></w:r>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:br/></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:br/></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'
)
$levantar.Range.InsertXML($levantarXml) | Out-Null

# ---------------------------------------------------------------------
# 2) PROJECT_VERSION=v1 -> add <w:lastRenderedPageBreak/>
# ---------------------------------------------------------------------
$projectVersion = Get-ParagraphByText("PROJECT_VERSION=v1")
$projectVersionXml = (
    '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>PROJECT_VERSION=v1</w:t></w:r>' +
    '</w:p>'
)
$projectVersion.Range.InsertXML($projectVersionXml) | Out-Null

# ---------------------------------------------------------------------
# 3) AWS_ACCESS_KEY_ID=... -> add <w:lastRenderedPageBreak/>
# ---------------------------------------------------------------------
$awsKey = Get-ParagraphByText("AWS_ACCESS_KEY_ID=AKIAXYKJTSSH3VRHA7HC")
$awsKeyXml = (
    '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>AWS_ACCESS_KEY_ID=AKIAXYKJTSSH3VRHA7HC</w:t></w:r>' +
    '</w:p>'
)
$awsKey.Range.InsertXML($awsKeyXml) | Out-Null

# ---------------------------------------------------------------------
# 4) Remove the old "_GoBack" bookmark after "... i dentro de la carpeta v1"
#    (it now lives on the new paragraph created in step 1).
# ---------------------------------------------------------------------
$npmPara = Get-ParagraphByText("i dentro de la carpeta v1")
$npmXml = (
    '<w:p><w:r><w:t xml:space="preserve">Comando </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> i dentro de la carpeta v1</w:t></w:r>' +
    '<w:r><w:br/></w:r>' +
    '</w:p>'
)
$npmPara.Range.InsertXML($npmXml) | Out-Null

Write-Output "Applied deploy-doc edits"
